# Generate Report for Handoff
# - Status flips from "Handed back: in sync with en-US" to "Ready for handoff"
#   on all three sheets (Overview, zh-cn, de-de).
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#   are refreshed to the new handoff-generation run.
# - The now-narrower status columns are resized to match.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-08-24 17:02:13"
$ws1.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws1.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---- zh-cn sheet ----
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("H2").Value = "2016-08-24 17:02:01"
$ws2.Columns.Item(3).ColumnWidth = 16.333333333333332

# ---- de-de sheet ----
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("H2").Value = "2016-08-24 17:02:13"
$ws3.Columns.Item(3).ColumnWidth = 16.333333333333332
